# Applies the "Add files via upload" edit to compliance_matrix.xlsx:
#  - E2:E29 used to show the literal text "Link" (hyperlinked). The author
#    overwrote that with the actual SharePoint evidence-folder URL text.
#  - Only E3 keeps a live hyperlink afterwards (now carrying a "display"
#    override equal to the new URL text); E2 and E4:E29 end up as plain
#    (non-hyperlinked) text.
#  - The worksheet view had scrolled right (topLeftCell=C1, selection G3);
#    the saved view now starts back at the default scroll with B10 selected.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newLinkText = "https://nflpk.sharepoint.com/nfl_site/IT%20Department/Forms/AllItems.aspx?id=%2Fnfl%5Fsite%2FIT%20Department%2FIT%2DGovernance%2FCOMPLIANCE%20MATRIX%2FEVIDENCE&viewid=15a9448d%2D5114%2D4f3f%2D81c3%2Ddaf0476b38d2&CT=1756293467065&OR=OWA%2DNT%2DMail&CID=bedfc76e%2D9e1e%2Db098%2D9316%2D9bdf1d6e274a&csf=1&web=1&e=svmXMI&FolderCTID=0x01200095CBBE0C68C0474BAF872F72573A9D97"
$e3Address = "https://nflpk.sharepoint.com/nfl_site/IT%20Department/Forms/AllItems.aspx?id=%2Fnfl%5Fsite%2FIT%20Department%2FIT%2DGovernance%2FCOMPLIANCE%20MATRIX%2FEVIDENCE%2FApplication%20Nfoods%2Ecom&viewid=15a9448d%2D5114%2D4f3f%2D81c3%2Ddaf0476b38d2&CT=1756293467065&OR=OWA%2DNT%2DMail&CID=bedfc76e%2D9e1e%2Db098%2D9316%2D9bdf1d6e274a&csf=1&web=1&e=svmXMI&FolderCTID=0x01200095CBBE0C68C0474BAF872F72573A9D97"

# All the other (mailto) hyperlinks on the sheet, in their original
# left-to-right-then-appended order, so they can be restored after the
# blanket Hyperlinks.Delete() below (this engine's Hyperlinks.Delete()
# removes every hyperlink on the sheet, not just the one it's called on).
# Split into the two runs that sandwich E3's hyperlink in the saved file.
$mailHyperlinksBeforeE3 = @(
    @("D3",  "mailto:farooquiyashal@gmail.com"),
    @("D2",  "mailto:farooquiyashal@gmail.com"),
    @("D4",  "mailto:farooquiyashal@gmail.com"),
    @("D5",  "mailto:farooquiyashal@gmail.com"),
    @("D6",  "mailto:farooquiyashal@gmail.com"),
    @("D7",  "mailto:aliyashal309@gmail.com"),
    @("D8",  "mailto:aliyashal309@gmail.com"),
    @("D9",  "mailto:aliyashal309@gmail.com"),
    @("D10", "mailto:aliyashal309@gmail.com"),
    @("D11", "mailto:aliyashal309@gmail.com"),
    @("D12", "mailto:aliyashal309@gmail.com"),
    @("D13", "mailto:aliyashal309@gmail.com"),
    @("D14", "mailto:aliyashal309@gmail.com"),
    @("D22", "mailto:yashal.ali@nfoods.com"),
    @("D23", "mailto:yashal.ali@nfoods.com"),
    @("D24", "mailto:yashal.ali@nfoods.com"),
    @("D25", "mailto:yashal.ali@nfoods.com"),
    @("D27", "mailto:aliyashal309@gmail.com"),
    @("D28", "mailto:aliyashal309@gmail.com"),
    @("D29", "mailto:aliyashal309@gmail.com")
)
$mailHyperlinksAfterE3 = @(
    @("D15", "mailto:farooquiyashal@gmail.com"),
    @("D16", "mailto:farooquiyashal@gmail.com"),
    @("D17", "mailto:farooquiyashal@gmail.com"),
    @("D18", "mailto:farooquiyashal@gmail.com"),
    @("D19", "mailto:farooquiyashal@gmail.com"),
    @("D20", "mailto:farooquiyashal@gmail.com"),
    @("D21", "mailto:farooquiyashal@gmail.com"),
    @("D26", "mailto:farooquiyashal@gmail.com")
)
$mailHyperlinks = $mailHyperlinksBeforeE3 + $mailHyperlinksAfterE3

# 1) Overwrite the displayed text of E2:E29 with the real URL (E3 is the
#    only one that keeps its hyperlink, re-added further down).
for ($row = 2; $row -le 29; $row++) {
    $ws.Range("E$row").Value2 = $newLinkText
}

# 2) This engine's Range.Hyperlinks.Delete() / Worksheet.Hyperlinks.Delete()
#    clears every hyperlink on the sheet in one go, so drop them all and
#    rebuild the set we want to keep (everything except the old E2 and
#    E4:E29 links).
$ws.Hyperlinks.Delete()

foreach ($pair in $mailHyperlinksBeforeE3) {
    $ws.Hyperlinks.Add($ws.Range($pair[0]), $pair[1]) | Out-Null
}

# E3 keeps its original (unchanged) target address but now also carries a
# "display" override equal to the new URL text typed into the cell.
# Add signature: Add(Anchor, Address, SubAddress, ScreenTip, TextToDisplay)
$ws.Hyperlinks.Add($ws.Range("E3"), $e3Address, "", "", $newLinkText) | Out-Null

foreach ($pair in $mailHyperlinksAfterE3) {
    $ws.Hyperlinks.Add($ws.Range($pair[0]), $pair[1]) | Out-Null
}

# Re-assert the "Hyperlink" cell style on every cell Hyperlinks.Add() just
# touched, and on the E column cells that lost their hyperlink, so nothing
# drifts from its original formatting (s="5" everywhere here).
foreach ($pair in $mailHyperlinks) {
    $ws.Range($pair[0]).Style = "Hyperlink"
}
for ($row = 2; $row -le 29; $row++) {
    $ws.Range("E$row").Style = "Hyperlink"
}

# 3) Restore the sheet view: scroll back to the default top-left and move
#    the selection to B10 (it had drifted to topLeftCell=C1 / G3 selected).
$ws.Range("B10").Select()
